# Applies the WALCL FRED data refresh:
#  - appends two new weekly observations to the "Data" sheet
#  - updates the refreshed-metadata fields on the "SeriesInfo" sheet

$wb = $excel.ActiveWorkbook

# --- Data sheet: append rows 98 and 99 ---------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Copy the formatting of the last existing data row (row 97) down into
# the two new rows so the date column keeps its date number format.
$wsData.Range("A97").Copy() | Out-Null
$wsData.Range("A98").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsData.Range("A99").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsData.Cells.Item(98, 1).Value = 45147
$wsData.Cells.Item(98, 2).Value = 8208.241

$wsData.Cells.Item(99, 1).Value = 45154
$wsData.Cells.Item(99, 2).Value = 8145.727

# --- SeriesInfo sheet: refresh metadata timestamps ----------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# B3/B4/B7 hold plain YYYY-MM-DD strings. Excel will happily reinterpret a
# literal "2023-08-22" typed into a General-formatted cell as a real date
# serial, so force text entry (format as Text, write, then drop back to the
# Normal style so the cell keeps its original, un-styled appearance).
$wsInfo.Range("B3").NumberFormat = "@"
$wsInfo.Range("B3").Value = "2023-08-22"
$wsInfo.Range("B3").Style = "Normal"

$wsInfo.Range("B4").NumberFormat = "@"
$wsInfo.Range("B4").Value = "2023-08-22"
$wsInfo.Range("B4").Style = "Normal"

$wsInfo.Range("B7").NumberFormat = "@"
$wsInfo.Range("B7").Value = "2023-08-16"
$wsInfo.Range("B7").Style = "Normal"

# B14 already contains a UTC-offset timestamp, which Excel's date parser
# won't accept as a date literal, so this one can be set directly.
$wsInfo.Range("B14").Value = "2023-08-17 15:33:36-05"
